$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right before
#    "dificuldade no gerenciamento de informações".
#    (_GoBack is a hidden bookmark, not shown by Bookmarks.Count /
#    enumeration, but it can still be looked up and deleted by name.)
# ------------------------------------------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 2) Replace the tail of the paragraph with the new wording.
# ------------------------------------------------------------------
$oldText = " à falta de automatização dos processos e transparência de dados na organização do negócio."
$newText = " à perda de dados e dificuldade de controlar as fichas do animal, receita, cobranças e dividas."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)

# ------------------------------------------------------------------
# 3) The diff splits the replaced text into several runs, and places
#    a (new) "_GoBack" bookmark between the "d" and the rest of
#    "dificuldade de controlar ...":
#
#      " à perda de dados e "   (run)
#      "d"                      (run)
#      <bookmarkStart _GoBack/><bookmarkEnd/>
#      "ificuldade de controlar as fichas do animal, receita, cobranças e dividas"  (run)
#      "."                      (run)
#
#    Inserting a zero-length bookmark at a text position forces Word
#    to split the run at that point, so we use that to create the
#    run boundaries, then remove the two throw-away helper bookmarks
#    while keeping the real "_GoBack" one.
# ------------------------------------------------------------------

$para = $d.Paragraphs.Item(4).Range
$text = $para.Text

$part1 = " à perda de dados e "
$part2 = "d"
$part3 = "ificuldade de controlar as fichas do animal, receita, cobranças e dividas"

$idx = $text.IndexOf($part1)
$p1End = $para.Start + $idx + $part1.Length
$p2End = $p1End + $part2.Length
$p3End = $p2End + $part3.Length

$d.Bookmarks.Add("ZZtempSplitA", $d.Range($p1End, $p1End)) | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($p2End, $p2End)) | Out-Null
$d.Bookmarks.Add("ZZtempSplitB", $d.Range($p3End, $p3End)) | Out-Null

$d.Bookmarks.Item("ZZtempSplitA").Delete()
$d.Bookmarks.Item("ZZtempSplitB").Delete()
